# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# The source data rows for a few matches had been built from the wrong
# fixture -> the row that was tagged with one match id actually carried
# the odds/result data belonging to a different match (and vice-versa).
# This fixes it by swapping/rotating the full data payload (every column
# except the running index in A and the constant Div/Date columns C/D)
# between the affected rows:
#   - rows 19 <-> 20
#   - rows 107 <-> 108
#   - rows 140 -> 141 -> 142 -> 140 (3-way rotation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns touched by the swap: B (id) and E..AD (HomeTeam .. PL_AhUnder).
# C (Div) and D (Date) are identical across the affected rows, so they are
# intentionally left untouched.
$colsToSwap = @(2) + @(5..30)

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-RowValues($ws, $row, $cols, $vals) {
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value = $vals[$c]
    }
}

# --- Swap rows 19 and 20 ---
$row19 = Get-RowValues $ws 19 $colsToSwap
$row20 = Get-RowValues $ws 20 $colsToSwap
Set-RowValues $ws 19 $colsToSwap $row20
Set-RowValues $ws 20 $colsToSwap $row19

# --- Swap rows 107 and 108 ---
$row107 = Get-RowValues $ws 107 $colsToSwap
$row108 = Get-RowValues $ws 108 $colsToSwap
Set-RowValues $ws 107 $colsToSwap $row108
Set-RowValues $ws 108 $colsToSwap $row107

# --- Rotate rows 140 -> 141 -> 142 -> 140 ---
# After the edit: row140 gets old row142 data, row141 gets old row140 data,
# row142 gets old row141 data.
$row140 = Get-RowValues $ws 140 $colsToSwap
$row141 = Get-RowValues $ws 141 $colsToSwap
$row142 = Get-RowValues $ws 142 $colsToSwap
Set-RowValues $ws 140 $colsToSwap $row142
Set-RowValues $ws 141 $colsToSwap $row140
Set-RowValues $ws 142 $colsToSwap $row141
